$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" header cells on the two existing sheets. ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last (3rd) tab. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header/date cell formats (bold+centered header, date number
# format on column A) instead of re-creating ad-hoc styles.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A17").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(45410.99999999999, 10, -0.8988364682179716, 20.81641198635106),
    @(45417.99999999999, 10, -0.8220813606906391, 20.23446357160307),
    @(45424.99999999999, 10, -1.153927426904526, 19.93840317422454),
    @(45445.99999999999, 9, -2.589071425492728, 19.69898527331735),
    @(45452.99999999999, 9, -1.449401050527757, 19.29020896049013),
    @(45466.99999999999, 9, -2.801745925535139, 19.28419490886729),
    @(45620.99999999999, 5, -5.493024570166915, 16.53868955929081),
    @(45634.99999999999, 5, -6.078325287617488, 15.72803567411644),
    @(45641.99999999999, 5, -6.060439587202504, 15.83773523691391),
    @(45648.99999999999, 5, -6.053303785942787, 14.9661093551234),
    @(45655.99999999999, 5, -6.10410879605315, 14.96791145835793),
    @(45662.99999999999, 4, -6.447804977349681, 14.51856768021199),
    @(45669.99999999999, 4, -6.738555475947982, 15.15736402992174),
    @(45676.99999999999, 4, -6.139035056484364, 13.67860852554622),
    @(45683.99999999999, 4, -6.750731875647188, 13.57136507048044),
    @(45690.99999999999, 4, -7.774466701624316, 15.15066115801454)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
